# Replace the old "Waktu Kampanye 2018 untuk Perseus: ..." line with the
# translated Gemini campaign dates, collapsing whatever runs made up that
# paragraph (the stray red space run, a leading "www.globeatnight.org"
# run, etc.) into a single plain run - matching the target OOXML which
# keeps the paragraph's <w:pPr> but has just one <w:r><w:t>...</w:t></w:r>.

$d = $word.ActiveDocument

$oldMarker  = "2018 untuk Perseus"
$newText    = "Waktu Kampanye Gemini: 14-23 Februari, 14-24 Maret"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $range = $para.Range

    if ($range.Text -notlike "*$oldMarker*") {
        continue
    }

    # Pull this paragraph's own pPr straight out of its current OOXML so
    # indentation/style/justification/etc. survive untouched.
    $wordOpenXml = $range.WordOpenXML
    $pPrXml = ""
    if ($wordOpenXml -match "(?s)<w:pPr>.*?</w:pPr>") {
        $pPrXml = $matches[0]
    }

    $newRunXml = "<w:r><w:t>" + $newText + "</w:t></w:r>"

    $packageXml = '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $pPrXml + $newRunXml + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $range.InsertXML($packageXml)
}
